$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price column (D) stores values as plain text even when they
# look numeric (e.g. "548.02"), so force text formatting before writing the
# updated prices to avoid Excel auto-converting them to numbers.
$ws.Range("D5,D6,D10,D12,D13,D14,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D31,D33,D34,D35,D36,D37,D39,D40,D41,D43,D44,D45,D46,D47,D49,D50").NumberFormat = "@"

$ws.Range('D2').Value = '61.440.11'
$ws.Range('E2').Value = '  -4.20%  '
$ws.Range('D3').Value = '3.005.92'
$ws.Range('E3').Value = '  -3.14%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '548.02'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '134.74'
$ws.Range('E6').Value = '  -4.68%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.002.79'
$ws.Range('E8').Value = '  -3.08%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '0.148'
$ws.Range('E10').Value = '  -5.89%  '
$ws.Range('E11').Value = '  -9.07%  '
$ws.Range('D12').Value = '0.450'
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '34.48'
$ws.Range('E13').Value = '  -1.64%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '0.0000221'
$ws.Range('E14').Value = '  -3.02%  '
$ws.Range('D15').Value = '3.494.71'
$ws.Range('E15').Value = '  -3.26%  '
$ws.Range('D16').Value = '61.534.63'
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').Value = '3.012.68'
$ws.Range('E18').Value = '  -3.08%  '
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = '471.82'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('D21').Value = '13.27'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').Value = '0.674'
$ws.Range('E22').Value = '  -4.18%  '
$ws.Range('D23').Value = '7.04'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('D24').Value = '80.00'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').Value = '12.08'
$ws.Range('E25').Value = '  -3.15%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = '2.71'
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('D28').Value = '7.81'
$ws.Range('E28').Value = '  -4.17%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('D31').Value = '25.66'
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('D33').Value = '5.52'
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('D34').Value = '2.30'
$ws.Range('E34').Value = '  -3.66%  '
$ws.Range('D35').Value = '55.44'
$ws.Range('E35').Value = '  -3.94%  '
$ws.Range('D36').Value = '5.89'
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('D37').Value = '453.59'
$ws.Range('E37').Value = '  -9.11%  '
$ws.Range('D38').Value = '3.187.96'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').Value = '0.0797'
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = '0.0382'
$ws.Range('E40').Value = '  -6.50%  '
$ws.Range('D41').Value = '0.117'
$ws.Range('E41').Value = '  -3.02%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').Value = '2.42'
$ws.Range('E43').Value = '  -12.19%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '26.12'
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('D46').Value = '0.244'
$ws.Range('E46').Value = '  -4.51%  '
$ws.Range('D47').Value = '1.98'
$ws.Range('E47').Value = '  -3.97%  '
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('D49').Value = '117.88'
$ws.Range('E49').Value = '  -4.97%  '
$ws.Range('B50').Value = 'BitgetToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D50').Value = '1.28'
$ws.Range('E50').Value = '  +6.75%  '
$ws.Range('B51').Value = 'PEPE'
$ws.Range('C51').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D51').Value = '0.0₃0491'
$ws.Range('E51').Value = '  -8.64%  '
